$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (index 0 / default) captured from an untouched cell,
# so we can force numeric-looking strings to stay text without altering
# the cell style/number-format of the workbook.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "56.605.04"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "2.490.39"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'492.18"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").Value = "'150.66"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +8.19%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").Value = "'0.516"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "2.497.27"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'5.77"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +5.80%  "
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").Value = "'0.336"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +2.72%  "
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").Value = "2.920.35"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "56.533.49"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "'21.20"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +3.50%  "
$ws.Range("D17").Value = "'0.0000137"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "2.493.13"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "'4.53"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("D20").Value = "'10.23"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").Value = "'321.37"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("E23").Value = "  +3.71%  "
$ws.Range("D24").Value = "'58.81"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").Value = "2.596.75"
$ws.Range("D29").Value = "'7.66"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +3.68%  "
$ws.Range("D30").Value = "0.0₃0799"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("D32").Value = "'150.46"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("D33").Value = "'18.37"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").Value = "'5.20"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").Value = "'1.17"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +5.14%  "
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("D38").Value = "'0.876"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +3.95%  "
$ws.Range("D39").Value = "'1.40"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +5.82%  "
$ws.Range("D40").Value = "'33.91"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("E41").Value = "  +3.48%  "
$ws.Range("D42").Value = "'0.0559"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("D43").Value = "'0.612"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "'0.993"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'4.85"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +8.53%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'264.86"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  +5.34%  "
$ws.Range("D47").Value = "'0.0927"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("D51").Value = "1.898.19"
$ws.Range("E51").Value = "  -3.83%  "
